# memory_grav.xlsx edit: append rows 40-46 to the "9.81" sheet with new
# grav-11 formula/value data, add a new "11" worksheet with its own 12 rows
# of data, and make the "9.81" sheet the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append rows 40-46 to the existing "9.81" sheet (below the existing
#    AVERAGE row 39, which is left untouched).
# ---------------------------------------------------------------------
$ws981 = $wb.Worksheets.Item("9.81")

$newRows981 = @(
    @('sub(sub(memProtectedDiv(cos(cos(sub(0, read(a0, protectedLog(a2))))), sub(read(a0, memProtectedDiv(a1, protectedLog(read(a0, protectedLog(0))))), a2)), a2), write(a0, write(a0, a1, 0), sin(sin(a2))))', -53.36, -57.29, -63.55, -101.51, -662.36, -445.79, -371.19, -376.79, -330.73, -625.26, -989.35, -897.12, -722.1, -1207.74, -1332.14, -1483.36),
    @('sub(read(a0, 0), add(add(add(add(write(a0, a1, add(memProtectedDiv(add(add(write(a0, read(a0, add(0, add(write(a0, a2, a2), read(a0, cos(a2))))), a1), a1), a1), conditional(a2, a1)), a1)), limit(a1, 0, a2)), a2), sin(read(a0, a1))), a2))', -114.66, -103.36, -111.13, -210.32, -279.37, -254.06, -250.68, -268.75, -269.74, -335.02, -371.18, -540.43, -712.85, -1411.94, -1595.65, -1710.47),
    @('sub(a1, add(memProtectedDiv(cos(0), limit(add(abs(a1), 0), 0, sub(sin(limit(a1, a2, memProtectedDiv(a2, sin(a1)))), read(a0, a2)))), write(a0, protectedLog(read(a0, write(a0, a1, protectedLog(a1)))), a2)))', -538.69, -584.94, -575.12, -669.16, -803.64, -712.44, -648.76, -547.39, -211.68, -224.8, -335.98, -1452.27, -1711.52, -1864.67, -1920.95, -2089.28),
    @('limit(conditional(a2, read(a0, a1)), write(a0, a2, a2), memProtectedDiv(memProtectedDiv(a2, cos(add(a1, a2))), protectedLog(sin(conditional(a1, memProtectedDiv(read(a0, limit(a1, a1, conditional(a1, memProtectedDiv(limit(abs(a1), 0, a2), memProtectedDiv(sub(read(a0, 0), limit(sub(a2, read(a0, read(a0, a2))), a1, a2)), abs(conditional(a2, 0))))))), a1))))))', -483.24, -500.67, -542.27, -631.86, -610.58, -613.95, -676.38, -794.32, -788.93, -977.56, -1490.9, -1650.19, -1812.14, -1919.24, -1994.35, -2055.18),
    @('sub(memProtectedDiv(sub(protectedLog(0), memProtectedDiv(memProtectedDiv(a1, sub(0, a1)), memProtectedDiv(sub(read(a0, a2), a2), cos(read(a0, a2))))), abs(write(a0, a1, sin(write(a0, 0, sin(a2)))))), a1)', -56.03, -55.73, -67.85, -114.04, -598.57, -565.79, -314.34, -260.4, -404.53, -682.04, -1085.91, -826.71, -750.12, -1210.49, -1341.26, -1476.33),
    @('sub(a1, memProtectedDiv(conditional(add(a2, read(a0, 0)), write(a0, write(a0, conditional(0, a2), 0), memProtectedDiv(add(a1, a1), a2))), a2))', -98.12, -128.26, -143.32, -158.86, -166.58, -172.25, -189.62, -200.35, -230.49, -293.33, -597.52, -932.99, -1051.29, -1260.92, -1583.19, -1751.11),
    @('conditional(read(a0, 0), write(a0, limit(a1, a1, limit(a2, a2, 0)), limit(a1, add(add(write(a0, a1, limit(a1, add(a2, memProtectedDiv(a1, a1)), a1)), memProtectedDiv(a1, a1)), cos(limit(a1, 0, sub(0, a2)))), 0)))', -580.48, -782.68, -833.86, -757.6, -523.56, -683.97, -655.31, -484.37, -523.18, -843.04, -1330.07, -1551.47, -1669.18, -1768.6, -1908.3, -1955.64)
)

$startRow981 = 40
for ($i = 0; $i -lt $newRows981.Count; $i++) {
    $r = $startRow981 + $i
    $rowdata = $newRows981[$i]
    for ($c = 1; $c -le $rowdata.Count; $c++) {
        $ws981.Cells.Item($r, $c).Value = $rowdata[$c - 1]
    }
}

# ---------------------------------------------------------------------
# 2. Add a brand new worksheet named "11" at the end of the workbook and
#    populate its 12 rows of formula/value data.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws11 = $wb.Worksheets.Add($null, $lastSheet)
$ws11.Name = "11"

$newRows11 = @(
    @('sub(add(a1, add(limit(a1, 0, read(a0, a1)), read(a0, a2))), write(a0, memProtectedDiv(a2, add(0, protectedLog(a1))), memProtectedDiv(sub(0, write(a0, 0, abs(write(a0, cos(abs(sub(a2, a1))), memProtectedDiv(a1, a2))))), add(sin(conditional(a2, write(a0, a2, add(a2, a2)))), limit(limit(0, memProtectedDiv(sin(a1), a1), sub(a1, a1)), 0, a2)))))', -357.43, -274.91, -238.23, -251.93, -331.23, -336.47, -283.45, -291.53, -283.39, -272.66, -342.42, -345.37, -1502.48, -1685.33, -1778.3, -1931.88),
    @('memProtectedDiv(sub(sin(read(a0, a1)), write(a0, 0, sin(a2))), abs(memProtectedDiv(sub(a2, sin(a2)), memProtectedDiv(sin(a1), read(a0, protectedLog(a1))))))', -238.18, -276.72, -319.23, -441.98, -596.57, -644.68, -591.77, -730.21, -755.57, -296.64, -382.86, -1297.22, -1489.24, -1680.67, -1829.78, -1931.78),
    @('sub(add(add(read(a0, protectedLog(abs(sin(a2)))), conditional(0, abs(a2))), sub(a2, memProtectedDiv(add(conditional(a1, 0), a1), a2))), write(a0, 0, sub(0, memProtectedDiv(add(write(a0, a1, conditional(memProtectedDiv(a1, a2), a1)), read(a0, a2)), a2))))', -852.38, -793.55, -931.95, -828.66, -935.81, -966.94, -827.02, -801.09, -341.52, -265.79, -380.61, -1503.4, -1720.24, -1843.7, -1917.37, -1979.78),
    @('limit(add(0, a2), conditional(a2, read(a0, limit(a2, a1, sin(a2)))), add(sub(a2, sub(memProtectedDiv(add(a1, read(a0, write(a0, 0, limit(a2, a1, a2)))), a2), sub(conditional(a1, a2), a1))), write(a0, memProtectedDiv(a2, 0), a2)))', -1641.73, -1811.87, -2144.95, -2090.29, -2016.75, -1949.89, -1856.26, -1739.01, -1517.65, -1073.6, -1566.47, -1762.24, -1887.93, -1989.55, -2057.86, -2126.05),
    @('sub(add(limit(a1, a2, a1), read(a0, add(a2, a2))), write(a0, sin(a1), add(memProtectedDiv(a1, conditional(a2, a1)), add(add(add(memProtectedDiv(0, 0), add(a2, a2)), add(0, a2)), write(a0, write(a0, a1, 0), a2)))))', -1879.12, -1989.46, -1880.64, -1897.77, -1848.04, -1739.85, -1674.16, -1477.89, -1336.82, -394.05, -1247.14, -1658.71, -1812.08, -1925.99, -2003.5, -2053.11),
    @('sub(read(a0, a1), write(a0, a1, memProtectedDiv(add(sub(limit(0, limit(a1, 0, memProtectedDiv(sub(a1, protectedLog(a1)), protectedLog(a1))), add(sub(a1, abs(sub(memProtectedDiv(a1, a2), a2))), protectedLog(limit(0, abs(sub(write(a0, a2, a2), add(sin(a2), 0))), a2)))), abs(sub(memProtectedDiv(a1, a2), read(a0, a1)))), protectedLog(limit(a2, 0, memProtectedDiv(add(sub(add(a1, a2), abs(sub(memProtectedDiv(a1, abs(sub(abs(0), add(sin(a2), 0)))), a2))), protectedLog(read(a0, read(a0, a2)))), a1)))), a2)))', -2015.15, -2062.02, -2087.44, -2132.41, -2071.12, -2049.1, -1959.36, -1820.44, -1544.35, -569.38, -1099.77, -1477.01, -1646.3, -1770.55, -1876.49, -1976.36),
    @('sub(memProtectedDiv(conditional(read(a0, cos(sub(0, a2))), write(a0, write(a0, add(a2, read(a0, 0)), protectedLog(write(a0, read(a0, 0), a1))), conditional(a2, conditional(read(a0, 0), cos(add(protectedLog(sub(a1, 0)), a1)))))), add(a2, memProtectedDiv(a2, 0))), sin(cos(a1)))', -1761.88, -1665.08, -1713.05, -1986.6, -2074.95, -1957.69, -1946.86, -1887.12, -1532.77, -596.8, -1065.85, -1663.95, -1784.88, -1872.64, -1959.12, -2045.08),
    @('sub(read(a0, a2), write(a0, a2, write(a0, write(a0, a2, write(a0, a1, add(write(a0, write(a0, a2, write(a0, a2, add(write(a0, protectedLog(a2), add(0, memProtectedDiv(memProtectedDiv(a2, protectedLog(cos(a2))), conditional(abs(a2), 0)))), memProtectedDiv(abs(a2), conditional(a1, 0))))), add(memProtectedDiv(a1, 0), memProtectedDiv(memProtectedDiv(a2, conditional(limit(a2, a2, 0), a2)), conditional(a2, 0)))), memProtectedDiv(sub(a2, memProtectedDiv(a2, a1)), conditional(a1, protectedLog(a1)))))), add(write(a0, a2, add(0, memProtectedDiv(memProtectedDiv(a2, protectedLog(cos(a2))), conditional(a2, limit(sin(a1), conditional(0, sub(cos(0), abs(a1))), sin(0)))))), memProtectedDiv(sub(a2, memProtectedDiv(a2, a1)), abs(a1))))))', -2049.74, -2268.23, -2175.2, -2164.05, -2106.42, -2140.33, -2005.4, -1777.37, -1529.51, -634.68, -1142.34, -1482.94, -1626.38, -1766.23, -1868.6, -1964.6),
    @('sub(sin(read(a0, 0)), limit(conditional(read(a0, a2), write(a0, cos(sub(a1, 0)), memProtectedDiv(cos(0), a2))), 0, read(a0, 0)))', -625.24, -694.03, -776.71, -886.64, -968.71, -1064.95, -1081.84, -1209.34, -1033.25, -417.1, -1390.01, -1567.33, -1632.22, -1771.4, -1923.03, -2006.75),
    @('sub(add(read(a0, 0), conditional(a1, add(a1, write(a0, limit(read(a0, abs(cos(a1))), abs(memProtectedDiv(limit(sin(a1), memProtectedDiv(a2, conditional(a1, write(a0, 0, read(a0, 0)))), sub(a1, 0)), sub(a1, sub(limit(sin(read(a0, cos(a1))), memProtectedDiv(a2, a1), sub(0, 0)), write(a0, a1, a1))))), limit(limit(memProtectedDiv(a1, a2), 0, abs(read(a0, a1))), 0, limit(limit(a1, 0, a1), abs(0), a2))), a2)))), write(a0, add(0, 0), memProtectedDiv(a1, conditional(a2, write(a0, a2, a1)))))', -579.83, -759.05, -724.01, -888.39, -893.73, -1043.51, -1019.99, -1061.83, -1000.21, -431.7, -1349.95, -1564.93, -1713.35, -1823.08, -1904.3, -2047.9),
    @('add(limit(a2, read(a0, a1), memProtectedDiv(memProtectedDiv(sub(a2, conditional(a1, sub(conditional(a2, cos(abs(a1))), protectedLog(sub(a2, a2))))), a1), a2)), cos(limit(cos(a1), a1, write(a0, conditional(0, a2), a2))))', -1838.53, -1986.86, -1820.35, -1859.26, -1729.15, -1689.91, -1554.81, -1563.12, -1360.68, -387.64, -1593.74, -1769.09, -1864.14, -1942.86, -1999.86, -2082.54),
    @('limit(read(a0, a1), a1, conditional(read(a0, a2), write(a0, a2, add(cos(a2), write(a0, a2, add(cos(sin(read(a0, protectedLog(abs(sub(a2, sin(a2))))))), add(a2, add(cos(0), a2))))))))', -1099.66, -1168.67, -1149.75, -1169.24, -1267.78, -1387.27, -1480.57, -1460.15, -1178.4, -555.6, -1223.32, -1442.37, -1673.58, -1782.93, -1874.49, -1958.91)
)

for ($i = 0; $i -lt $newRows11.Count; $i++) {
    $r = $i + 1
    $rowdata = $newRows11[$i]
    for ($c = 1; $c -le $rowdata.Count; $c++) {
        $ws11.Cells.Item($r, $c).Value = $rowdata[$c - 1]
    }
}

# ---------------------------------------------------------------------
# 3. Make the "9.81" sheet the active tab (activeTab goes from 3 to 0).
# ---------------------------------------------------------------------
$ws981.Activate()
